$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by
# Excels type inference (single-decimal numeric-looking strings such as
# "1.004"). We briefly mark them as Text, assign the literal string, then
# restore the Normal style so no stray number-format/style id lingers on
# the cell (keeps the saved XML free of spurious s="" attributes).
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '27.969.31'
$ws.Range('E2').Value = '  +1.86%  '
$ws.Range('D3').Value = '1.904.63'
$ws.Range('E3').Value = '  +2.32%  '
Set-TextValue 'D4' '1.004'
$ws.Range('E4').Value = '  -0.68%  '
Set-TextValue 'D5' '317.41'
$ws.Range('E5').Value = '  +1.96%  '
$ws.Range('E6').Value = '  -0.72%  '
Set-TextValue 'D7' '0.4818'
$ws.Range('E7').Value = '  +0.99%  '
Set-TextValue 'D8' '0.3804'
$ws.Range('E8').Value = '  -0.40%  '
Set-TextValue 'D9' '0.07362'
$ws.Range('E9').Value = '  +0.55%  '
Set-TextValue 'D10' '0.9309'
$ws.Range('E10').Value = '  -0.05%  '
Set-TextValue 'D11' '20.80'
$ws.Range('E11').Value = '  +0.05%  '
Set-TextValue 'D12' '0.07744'
$ws.Range('E12').Value = '  -0.58%  '
$ws.Range('D13').Value = '1.888.36'
$ws.Range('E13').Value = '  +1.04%  '
Set-TextValue 'D14' '5.485'
$ws.Range('E14').Value = '  +0.69%  '
Set-TextValue 'D15' '6.639'
$ws.Range('E15').Value = '  +1.27%  '
Set-TextValue 'D16' '91.59'
$ws.Range('E16').Value = '  +1.59%  '
Set-TextValue 'D17' '1.004'
$ws.Range('E17').Value = '  -0.74%  '
Set-TextValue 'D18' '0.000008892'
$ws.Range('E18').Value = '  +0.76%  '
Set-TextValue 'D19' '1.002'
$ws.Range('E19').Value = '  -0.68%  '
$ws.Range('D20').Value = '28.003.66'
$ws.Range('E20').Value = '  +1.88%  '
Set-TextValue 'D21' '14.68'
$ws.Range('E21').Value = '  +0.46%  '
Set-TextValue 'D22' '5.138'
$ws.Range('E22').Value = '  +0.81%  '
$ws.Range('D23').Value = '2.153.52'
$ws.Range('E23').Value = '  +2.31%  '
Set-TextValue 'D24' '10.88'
$ws.Range('E24').Value = '  +1.76%  '
Set-TextValue 'D25' '156.12'
$ws.Range('E25').Value = '  +0.50%  '
Set-TextValue 'D26' '1.911'
$ws.Range('E26').Value = '  -1.53%  '
Set-TextValue 'D27' '18.50'
$ws.Range('E27').Value = '  +0.16%  '
Set-TextValue 'D28' '2.110'
$ws.Range('E28').Value = '  +4.88%  '
Set-TextValue 'D29' '117.40'
$ws.Range('E29').Value = '  +1.70%  '
Set-TextValue 'D30' '4.967'
$ws.Range('E30').Value = '  +0.23%  '
Set-TextValue 'D31' '0.08944'
$ws.Range('E31').Value = '  +0.41%  '
Set-TextValue 'D32' '3.249'
$ws.Range('E32').Value = '  -2.36%  '
Set-TextValue 'D33' '1.251'
$ws.Range('E33').Value = '  +3.69%  '
Set-TextValue 'D34' '0.7710'
$ws.Range('E34').Value = '  +2.37%  '
Set-TextValue 'D35' '4.661'
$ws.Range('E35').Value = '  +1.39%  '
Set-TextValue 'D36' '2.595'
$ws.Range('E36').Value = '  -4.29%  '
Set-TextValue 'D37' '0.02046'
$ws.Range('E37').Value = '  +0.14%  '
Set-TextValue 'D38' '1.105'
$ws.Range('E38').Value = '  -1.85%  '
Set-TextValue 'D39' '0.5506'
$ws.Range('E39').Value = '  -1.03%  '
Set-TextValue 'D40' '0.05283'
$ws.Range('E40').Value = '  +0.25%  '
Set-TextValue 'D41' '2.995'
$ws.Range('E41').Value = '  +0.23%  '
Set-TextValue 'D42' '6.978'
$ws.Range('E42').Value = '  -0.89%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D43' '0.1529'
$ws.Range('E43').Value = '  +0.40%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D44' '8.472'
$ws.Range('E44').Value = '  -1.84%  '
Set-TextValue 'D45' '110.12'
$ws.Range('E45').Value = '  +6.93%  '
Set-TextValue 'D46' '10.70'
$ws.Range('E46').Value = '  -0.13%  '
Set-TextValue 'D47' '0.4820'
$ws.Range('E47').Value = '  -1.32%  '
Set-TextValue 'D48' '1.002'
$ws.Range('E48').Value = '  -0.75%  '
Set-TextValue 'D49' '1.643'
$ws.Range('E49').Value = '  -1.44%  '
Set-TextValue 'D50' '67.87'
$ws.Range('E50').Value = '  +0.65%  '
Set-TextValue 'D51' '0.06077'
$ws.Range('E51').Value = '  -0.24%  '
